$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Capture the existing header-row date labels BEFORE we shift anything, so we
# can relocate them to their new positions once the two new columns exist.
# Before:  B1 = "Jun_13"   C1 = "Jun_10"
# After :  B1 = "Jun_17"   C1 = "Jun_15"   D1 = "Jun_13"   E1 = "Jun_10"
# ---------------------------------------------------------------------------
$oldB1 = $ws.Range("B1").Value2

# ---------------------------------------------------------------------------
# Insert two new (currently blank) columns at C:D. This pushes the former
# column C (rating-change detail / style) two places to the right, landing
# in column E, while leaving columns A and B untouched for the data rows.
# ---------------------------------------------------------------------------
$ws.Columns("C:D").Insert()

# ---------------------------------------------------------------------------
# Header row: the value that used to live in B1 now belongs in D1, and the
# two freshly inserted header cells get the two new snapshot-date labels.
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = $oldB1

$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# ---------------------------------------------------------------------------
# Data rows (2-27): the two new columns are placeholders carrying the same
# "UN" marker used throughout column B, styled like the rest of the table.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# ---------------------------------------------------------------------------
# Column widths: columns C and D mirror the original column C's width, and
# the shifted-out column E keeps that same width too.
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14
$ws.Columns("E").ColumnWidth = 7.14
